$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cell C10 ("From" value of rule R30) changes from 18 to 1
$ws.Range("C10").Value = 1
